$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the entire row 658 ("「うちゅうへとびたい」..." post), which shifts
# all subsequent rows (659-859) up by one, turning the sheet's extent from
# A1:C859 into A1:C858.
$ws.Rows.Item(658).Delete()
